# daily auto push: 2026-01-09 06:50 UTC
# A new observation row for 2026/01/09 (Fri) 13:00, ranking 23 was inserted
# right before the existing 2026/12/29 block (sheet row 586), pushing every
# row from the old 586 down through the old 627 (last row) one position down
# (new rows 587-628). The used range grows from A1:D627 to A1:D628.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at sheet row 586 - shifts rows 586..627 down to 587..628
$ws.Rows.Item(586).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not real
# Excel dates, in this workbook. Force text interpretation before assigning
# so "2026/01/09" isn't auto-converted into a date serial number, then drop
# back to the default "Normal" style so the cell carries no style index -
# matching every other data row in the sheet.
$ws.Range("A586").NumberFormat = "@"
$ws.Range("A586").Value = "2026/01/09"
$ws.Range("A586").Style = "Normal"

$ws.Range("B586").Value = "金"
$ws.Range("C586").Value = 13
$ws.Range("D586").Value = 23
